$d = $word.ActiveDocument

# --- 1) Swap the font used throughout the document: Zapfino -> Gill Sans Light ---
# Update every run's complex-script / East Asian / ASCII / "other" font names so
# that w:rFonts ends up with ascii/cs/eastAsia/hAnsi all set to "Gill Sans Light".
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark so pPr/rPr isn't touched
    if ($r.Font.NameAscii -eq "Zapfino") {
        $r.Font.NameAscii   = "Gill Sans Light"
        $r.Font.NameFarEast = "Gill Sans Light"
        $r.Font.NameOther   = "Gill Sans Light"
        $r.Font.NameBi      = "Gill Sans Light"
    }
}

# --- 2) Update the wording of each paragraph's text ---
$replacements = @(
    @{ Old = "This text is written in Zapfino and is using no ligatures.";
       New = "This field test text is written in Gill Sans Light and is using no ligatures." },
    @{ Old = "This text is written in Zapfino and is using default ligatures.";
       New = "This field test text is written in Gill Sans Light and is using default ligatures." },
    @{ Old = "This text is written in Zapfino and is using all supported ligatures.";
       New = "This field test text is written in Gill Sans Light and is using all supported ligatures." },
    @{ Old = "The attributes of this text do not mention ligatures, so the text should be using default ligatures. It is written in Zapfino, by the way.";
       New = "The attributes of this field test text do not mention ligatures, so the text should be using default ligatures. It is written in Gill Sans Light, by the way." }
)

foreach ($rep in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($rep.Old, $true, $true, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null
}

Write-Output "edit complete"
